$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "59-2="
$t.Cell(1,2).Range.Text = "72+18="
$t.Cell(1,3).Range.Text = "57-49="
$t.Cell(1,4).Range.Text = "73-25="
$t.Cell(1,5).Range.Text = "47-30="

$t.Cell(2,1).Range.Text = "64-21="
$t.Cell(2,2).Range.Text = "62-53="
$t.Cell(2,3).Range.Text = "73+9="
$t.Cell(2,4).Range.Text = "98-15="
$t.Cell(2,5).Range.Text = "11+25="

$t.Cell(3,1).Range.Text = "73+18="
$t.Cell(3,2).Range.Text = "39-38="
$t.Cell(3,3).Range.Text = "59-3="
$t.Cell(3,4).Range.Text = "74-56="
$t.Cell(3,5).Range.Text = "9+24="

$t.Cell(4,1).Range.Text = "71-10="
$t.Cell(4,2).Range.Text = "21+11="
$t.Cell(4,3).Range.Text = "31+16="
$t.Cell(4,4).Range.Text = "51-31="
$t.Cell(4,5).Range.Text = "8+78="

$t.Cell(5,1).Range.Text = "7+61="
$t.Cell(5,2).Range.Text = "52-4="
$t.Cell(5,3).Range.Text = "33-18="
$t.Cell(5,4).Range.Text = "39+20="
$t.Cell(5,5).Range.Text = "74-69="

$t.Cell(6,1).Range.Text = "54+13="
$t.Cell(6,2).Range.Text = "21+68="
$t.Cell(6,3).Range.Text = "35+44="
$t.Cell(6,4).Range.Text = "60-13="
$t.Cell(6,5).Range.Text = "50-42="

$t.Cell(7,1).Range.Text = "9+29="
$t.Cell(7,2).Range.Text = "13+8="
$t.Cell(7,3).Range.Text = "37-22="
$t.Cell(7,4).Range.Text = "54+5="
$t.Cell(7,5).Range.Text = "92-36="

$t.Cell(8,1).Range.Text = "78+3="
$t.Cell(8,2).Range.Text = "77-3="
$t.Cell(8,3).Range.Text = "17+8="
$t.Cell(8,4).Range.Text = "40-23="
$t.Cell(8,5).Range.Text = "89-63="

$t.Cell(9,1).Range.Text = "9+26="
$t.Cell(9,2).Range.Text = "68+18="
$t.Cell(9,3).Range.Text = "65+11="
$t.Cell(9,4).Range.Text = "31+41="
$t.Cell(9,5).Range.Text = "27+69="

$t.Cell(10,1).Range.Text = "26+5="
$t.Cell(10,2).Range.Text = "18-11="
$t.Cell(10,3).Range.Text = "6+71="
$t.Cell(10,4).Range.Text = "41-22="
$t.Cell(10,5).Range.Text = "86+8="

$t.Cell(11,1).Range.Text = "68-0="
$t.Cell(11,2).Range.Text = "50-3="
$t.Cell(11,3).Range.Text = "6+1="
$t.Cell(11,4).Range.Text = "26+16="
$t.Cell(11,5).Range.Text = "32+42="

$t.Cell(12,1).Range.Text = "11-3="
$t.Cell(12,2).Range.Text = "30+33="
$t.Cell(12,3).Range.Text = "8+64="
$t.Cell(12,4).Range.Text = "11+5="
$t.Cell(12,5).Range.Text = "14+27="

$t.Cell(13,1).Range.Text = "80+7="
$t.Cell(13,2).Range.Text = "39-37="
$t.Cell(13,3).Range.Text = "14+61="
$t.Cell(13,4).Range.Text = "30-18="
$t.Cell(13,5).Range.Text = "45-26="

$t.Cell(14,1).Range.Text = "39-28="
$t.Cell(14,2).Range.Text = "80-3="
$t.Cell(14,3).Range.Text = "41-31="
$t.Cell(14,4).Range.Text = "3+20="
$t.Cell(14,5).Range.Text = "32+39="

$t.Cell(15,1).Range.Text = "28+4="
$t.Cell(15,2).Range.Text = "47+24="
$t.Cell(15,3).Range.Text = "75-61="
$t.Cell(15,4).Range.Text = "40+49="
$t.Cell(15,5).Range.Text = "91-26="

$t.Cell(16,1).Range.Text = "2+90="
$t.Cell(16,2).Range.Text = "98-80="
$t.Cell(16,3).Range.Text = "11+44="
$t.Cell(16,4).Range.Text = "92-24="
$t.Cell(16,5).Range.Text = "65-40="

$t.Cell(17,1).Range.Text = "1+94="
$t.Cell(17,2).Range.Text = "87+7="
$t.Cell(17,3).Range.Text = "38-30="
$t.Cell(17,4).Range.Text = "97-86="
$t.Cell(17,5).Range.Text = "92-37="

$t.Cell(18,1).Range.Text = "21+22="
$t.Cell(18,2).Range.Text = "28+15="
$t.Cell(18,3).Range.Text = "99-51="
$t.Cell(18,4).Range.Text = "40-37="
$t.Cell(18,5).Range.Text = "18-15="

$t.Cell(19,1).Range.Text = "54+8="
$t.Cell(19,2).Range.Text = "28-20="
$t.Cell(19,3).Range.Text = "68-59="
$t.Cell(19,4).Range.Text = "36+6="
$t.Cell(19,5).Range.Text = "65+3="

$t.Cell(20,1).Range.Text = "68+10="
$t.Cell(20,2).Range.Text = "17+56="
$t.Cell(20,3).Range.Text = "14+21="
$t.Cell(20,4).Range.Text = "29+26="
$t.Cell(20,5).Range.Text = "56+11="
